# Update the "last updated" timestamp banner (sharedStrings text change)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 16 de Agosto de 2020 a las 04:36"

# Refresh per-country Covid figures. A handful of countries (Australia/Austria,
# Paraguay/Zambia/Malasia, Libia/Finlandia/Luxemburgo and the four-way cluster around
# Mongolia) swap rank with their neighbours so the table stays sorted by "Casos totales" desc.

$ws.Cells.Item(31, 1).Value = "Bolivia"
$ws.Cells.Item(31, 2).Value = 99146
$ws.Cells.Item(31, 3).Value = 1196
$ws.Cells.Item(31, 4).Value = 35638
$ws.Cells.Item(31, 5).Value = 59505
$ws.Cells.Item(31, 6).Value = 0
$ws.Cells.Item(31, 7).Value = 64
$ws.Cells.Item(31, 8).Value = 4003

$ws.Cells.Item(40, 1).Value = "Belgica"
$ws.Cells.Item(40, 2).Value = 77869
$ws.Cells.Item(40, 3).Value = 756
$ws.Cells.Item(40, 4).Value = 17981
$ws.Cells.Item(40, 5).Value = 49953
$ws.Cells.Item(40, 6).Value = 0
$ws.Cells.Item(40, 7).Value = 11
$ws.Cells.Item(40, 8).Value = 9935

$ws.Cells.Item(51, 1).Value = "Honduras"
$ws.Cells.Item(51, 2).Value = 49979
$ws.Cells.Item(51, 3).Value = 512
$ws.Cells.Item(51, 4).Value = 7255
$ws.Cells.Item(51, 5).Value = 41157
$ws.Cells.Item(51, 6).Value = 0
$ws.Cells.Item(51, 7).Value = 19
$ws.Cells.Item(51, 8).Value = 1567

$ws.Cells.Item(71, 1).Value = "Australia"
$ws.Cells.Item(71, 2).Value = 23287
$ws.Cells.Item(71, 3).Value = 252
$ws.Cells.Item(71, 4).Value = 13634
$ws.Cells.Item(71, 5).Value = 9257
$ws.Cells.Item(71, 6).Value = 0
$ws.Cells.Item(71, 7).Value = 17
$ws.Cells.Item(71, 8).Value = 396

$ws.Cells.Item(72, 1).Value = "Austria"
$ws.Cells.Item(72, 2).Value = 23179
$ws.Cells.Item(72, 3).Value = 0
$ws.Cells.Item(72, 4).Value = 20627
$ws.Cells.Item(72, 5).Value = 1824
$ws.Cells.Item(72, 6).Value = 0
$ws.Cells.Item(72, 7).Value = 0
$ws.Cells.Item(72, 8).Value = 728

$ws.Cells.Item(80, 1).Value = "Corea del Sur"
$ws.Cells.Item(80, 2).Value = 15318
$ws.Cells.Item(80, 3).Value = 279
$ws.Cells.Item(80, 4).Value = 13910
$ws.Cells.Item(80, 5).Value = 1103
$ws.Cells.Item(80, 6).Value = 0
$ws.Cells.Item(80, 7).Value = 0
$ws.Cells.Item(80, 8).Value = 305

$ws.Cells.Item(88, 1).Value = "Paraguay"
$ws.Cells.Item(88, 2).Value = 9381
$ws.Cells.Item(88, 3).Value = 0
$ws.Cells.Item(88, 4).Value = 5841
$ws.Cells.Item(88, 5).Value = 3413
$ws.Cells.Item(88, 6).Value = 0
$ws.Cells.Item(88, 7).Value = 0
$ws.Cells.Item(88, 8).Value = 127

$ws.Cells.Item(89, 1).Value = "Zambia"
$ws.Cells.Item(89, 2).Value = 9186
$ws.Cells.Item(89, 3).Value = 0
$ws.Cells.Item(89, 4).Value = 8065
$ws.Cells.Item(89, 5).Value = 861
$ws.Cells.Item(89, 6).Value = 0
$ws.Cells.Item(89, 7).Value = 0
$ws.Cells.Item(89, 8).Value = 260

$ws.Cells.Item(90, 1).Value = "Malasia"
$ws.Cells.Item(90, 2).Value = 9175
$ws.Cells.Item(90, 3).Value = 0
$ws.Cells.Item(90, 4).Value = 8831
$ws.Cells.Item(90, 5).Value = 219
$ws.Cells.Item(90, 6).Value = 0
$ws.Cells.Item(90, 7).Value = 0
$ws.Cells.Item(90, 8).Value = 125

$ws.Cells.Item(97, 1).Value = "Libia"
$ws.Cells.Item(97, 2).Value = 7738
$ws.Cells.Item(97, 3).Value = 0
$ws.Cells.Item(97, 4).Value = 894
$ws.Cells.Item(97, 5).Value = 6699
$ws.Cells.Item(97, 6).Value = 0
$ws.Cells.Item(97, 7).Value = 0
$ws.Cells.Item(97, 8).Value = 145

$ws.Cells.Item(98, 1).Value = "Finlandia"
$ws.Cells.Item(98, 2).Value = 7720
$ws.Cells.Item(98, 3).Value = 0
$ws.Cells.Item(98, 4).Value = 7050
$ws.Cells.Item(98, 5).Value = 337
$ws.Cells.Item(98, 6).Value = 0
$ws.Cells.Item(98, 7).Value = 0
$ws.Cells.Item(98, 8).Value = 333

$ws.Cells.Item(99, 1).Value = "Luxemburgo"
$ws.Cells.Item(99, 2).Value = 7439
$ws.Cells.Item(99, 3).Value = 0
$ws.Cells.Item(99, 4).Value = 6500
$ws.Cells.Item(99, 5).Value = 816
$ws.Cells.Item(99, 6).Value = 0
$ws.Cells.Item(99, 7).Value = 0
$ws.Cells.Item(99, 8).Value = 123

$ws.Cells.Item(140, 1).Value = "Nueva Zelanda"
$ws.Cells.Item(140, 2).Value = 1622
$ws.Cells.Item(140, 3).Value = 13
$ws.Cells.Item(140, 4).Value = 1531
$ws.Cells.Item(140, 5).Value = 69
$ws.Cells.Item(140, 6).Value = 0
$ws.Cells.Item(140, 7).Value = 0
$ws.Cells.Item(140, 8).Value = 22

$ws.Cells.Item(178, 1).Value = "San Martin (Parte Holandesa)"
$ws.Cells.Item(178, 2).Value = 300
$ws.Cells.Item(178, 3).Value = 31
$ws.Cells.Item(178, 4).Value = 107
$ws.Cells.Item(178, 5).Value = 176
$ws.Cells.Item(178, 6).Value = 0
$ws.Cells.Item(178, 7).Value = 0
$ws.Cells.Item(178, 8).Value = 17

$ws.Cells.Item(179, 1).Value = "Islas Turcas y Caicos"
$ws.Cells.Item(179, 2).Value = 298
$ws.Cells.Item(179, 3).Value = 24
$ws.Cells.Item(179, 4).Value = 55
$ws.Cells.Item(179, 5).Value = 241
$ws.Cells.Item(179, 6).Value = 0
$ws.Cells.Item(179, 7).Value = 0
$ws.Cells.Item(179, 8).Value = 2

$ws.Cells.Item(180, 1).Value = "Mongolia"
$ws.Cells.Item(180, 2).Value = 298
$ws.Cells.Item(180, 3).Value = 0
$ws.Cells.Item(180, 4).Value = 272
$ws.Cells.Item(180, 5).Value = 26
$ws.Cells.Item(180, 6).Value = 0
$ws.Cells.Item(180, 7).Value = 0
$ws.Cells.Item(180, 8).Value = 0

$ws.Cells.Item(181, 1).Value = "Eritrea"
$ws.Cells.Item(181, 2).Value = 285
$ws.Cells.Item(181, 3).Value = 0
$ws.Cells.Item(181, 4).Value = 248
$ws.Cells.Item(181, 5).Value = 37
$ws.Cells.Item(181, 6).Value = 0
$ws.Cells.Item(181, 7).Value = 0
$ws.Cells.Item(181, 8).Value = 0

$ws.Cells.Item(182, 1).Value = "Camboya"
$ws.Cells.Item(182, 2).Value = 273
$ws.Cells.Item(182, 3).Value = 0
$ws.Cells.Item(182, 4).Value = 238
$ws.Cells.Item(182, 5).Value = 35
$ws.Cells.Item(182, 6).Value = 0
$ws.Cells.Item(182, 7).Value = 0
$ws.Cells.Item(182, 8).Value = 0

$ws.Cells.Item(183, 1).Value = "Papua Nueva Guinea"
$ws.Cells.Item(183, 2).Value = 271
$ws.Cells.Item(183, 3).Value = 0
$ws.Cells.Item(183, 4).Value = 78
$ws.Cells.Item(183, 5).Value = 190
$ws.Cells.Item(183, 6).Value = 0
$ws.Cells.Item(183, 7).Value = 0
$ws.Cells.Item(183, 8).Value = 3
